$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(5, 1).Value = 42604.891412037039
$ws.Cells.Item(5, 2).Value = "Noun"
$ws.Cells.Item(5, 3).Value = 13111
$ws.Cells.Item(5, 4).Value = 8376
$ws.Cells.Item(5, 5).Value = 1368
$ws.Cells.Item(5, 6).Value = 191
$ws.Cells.Item(5, 7).Value = 140
$ws.Cells.Item(5, 8).Value = 57
$ws.Cells.Item(5, 9).Value = 42
$ws.Cells.Item(5, 10).Value = 4
$ws.Cells.Item(5, 11).Value = 20
$ws.Cells.Item(5, 12).Value = 16
$ws.Cells.Item(5, 13).Value = 83
